$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Highlight the existing row 7 (Change_ID=7) accuracy values in yellow (D14:G14, D15:G15)
$ws.Range("D14:G15").Interior.Color = 65535

# Fill in the accuracy values for row 17 (Change_ID 9 / scheduler) and highlight the row
$ws.Range("D17").Value = 70.040000000000006
$ws.Range("E17").Value = 43.51
$ws.Range("F17").Value = 89.94
$ws.Range("G17").Value = 88.33
$ws.Range("A17:G17").Interior.Color = 65535

# Add a new row for Change_ID 10: final_div_factor set to None
$ws.Range("A18").Value = 10
$ws.Range("B18").Value = "final_div_factor"
$ws.Range("C18").Value = "None"

# Update the active selection to match the recorded state
$ws.Range("D18").Select()
